# Trade #22 closed at 2026-02-17 15:18:59 - unknown UNKNOWN +0.000%
#
# Updates the aggregated summary / status figures to reflect the newly
# closed trade, and appends the new trade row (#22, 0-indexed "22" in
# the Trade # column) to both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.76   # Current Capital
$summary.Range("B5").Value = -0.22     # Total P&L %
$summary.Range("B6").Value = 22        # Total Trades
$summary.Range("B9").Value = 22.73     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.76000000000001   # Capital
$status.Range("D4").Value = 22                  # Trades
$status.Range("F4").Value = -0.24               # P&L %
$status.Range("G4").Value = 22.73               # Win Rate %

# ---------------------------------------------------------------------
# Helper that appends the new trade row (row 23) to a trades sheet
# ---------------------------------------------------------------------
function Add-Trade23Row($ws) {
    $ws.Cells.Item(23, 1).Value = 22

    $ws.Cells.Item(23, 2).NumberFormat = "@"
    $ws.Cells.Item(23, 2).Value = "2026-02-17"
    $ws.Cells.Item(23, 2).Style = "Normal"

    $ws.Cells.Item(23, 3).NumberFormat = "@"
    $ws.Cells.Item(23, 3).Value = "15:18:53"
    $ws.Cells.Item(23, 3).Style = "Normal"

    $ws.Cells.Item(23, 4).Value = "MarketMaking"
    $ws.Cells.Item(23, 5).Value = "UP"
    $ws.Cells.Item(23, 6).Value = 0.8649829999999999
    $ws.Cells.Item(23, 7).Value = 0.86
    $ws.Cells.Item(23, 8).Value = "CLOSED"
    $ws.Cells.Item(23, 9).Value = -0.5760999999999999
    $ws.Cells.Item(23, 10).Value = -0
    $ws.Cells.Item(23, 11).Value = 99.76000000000001
    $ws.Cells.Item(23, 12).Value = 0
    $ws.Cells.Item(23, 13).Value = 0
    $ws.Cells.Item(23, 14).Value = 0.6
    $ws.Cells.Item(23, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(23, 16).Value = "early_exit"
    $ws.Cells.Item(23, 17).Value = 0.14
}

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade23Row $allTrades

# ---------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade23Row $marketMaking
